$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08097566666666667
$ws.Range("H2").Value = 0.242927
$ws.Range("I2").Value = 0.005588990034505014
$ws.Range("J2").Value = 0.005588990034505015
$ws.Range("M2").Value = 14.11187666666667
$ws.Range("N2").Value = 42.33562999999999
$ws.Range("O2").Value = 0.08862966207485527
$ws.Range("P2").Value = 0.08862966207485526
$ws.Range("Q2").Value = 1.142718621001111
$ws.Range("R2").Value = 10.28446758901
$ws.Range("S2").Value = 0.0004953502980979131
$ws.Range("T2").Value = 0.0004953502980979131

$ws.Range("G3").Value = 0.08097566666666667
$ws.Range("H3").Value = 0.242927
$ws.Range("I3").Value = 0.005588990034505014
$ws.Range("J3").Value = 0.005588990034505015
$ws.Range("O3").Value = 0.7176943460983047
$ws.Range("P3").Value = 0.7176943460983046
$ws.Range("Q3").Value = 9.253365907917889
$ws.Range("R3").Value = 83.28029317126101
$ws.Range("S3").Value = 0.004011186548164018
$ws.Range("T3").Value = 0.004011186548164018

$ws.Range("G4").Value = 0.08097566666666667
$ws.Range("H4").Value = 0.242927
$ws.Range("I4").Value = 0.005588990034505014
$ws.Range("J4").Value = 0.005588990034505015
$ws.Range("O4").Value = 0.1936759918268401
$ws.Range("P4").Value = 0.1936759918268401
$ws.Range("Q4").Value = 2.497100373850778
$ws.Range("R4").Value = 22.473903364657
$ws.Range("S4").Value = 0.001082453188243084
$ws.Range("T4").Value = 0.001082453188243084

$ws.Range("I5").Value = 0.6976944377922635
$ws.Range("J5").Value = 0.6976944377922635
$ws.Range("M5").Value = 14.11187666666667
$ws.Range("N5").Value = 42.33562999999999
$ws.Range("O5").Value = 0.08862966207485527
$ws.Range("P5").Value = 0.08862966207485526
$ws.Range("Q5").Value = 142.64982061374
$ws.Range("R5").Value = 1283.84838552366
$ws.Range("S5").Value = 0.06183642225303444
$ws.Range("T5").Value = 0.06183642225303444

$ws.Range("I6").Value = 0.6976944377922635
$ws.Range("J6").Value = 0.6976944377922635
$ws.Range("O6").Value = 0.7176943460983047
$ws.Range("P6").Value = 0.7176943460983046
$ws.Range("S6").Value = 0.5007313533077429
$ws.Range("T6").Value = 0.5007313533077428

$ws.Range("I7").Value = 0.6976944377922635
$ws.Range("J7").Value = 0.6976944377922635
$ws.Range("O7").Value = 0.1936759918268401
$ws.Range("P7").Value = 0.1936759918268401
$ws.Range("S7").Value = 0.1351266622314863
$ws.Range("T7").Value = 0.1351266622314862

$ws.Range("I8").Value = 0.2967165721732315
$ws.Range("J8").Value = 0.2967165721732316
$ws.Range("M8").Value = 14.11187666666667
$ws.Range("N8").Value = 42.33562999999999
$ws.Range("O8").Value = 0.08862966207485527
$ws.Range("P8").Value = 0.08862966207485526
$ws.Range("Q8").Value = 60.66633686742666
$ws.Range("R8").Value = 545.99703180684
$ws.Range("S8").Value = 0.02629788952372291
$ws.Range("T8").Value = 0.02629788952372291

$ws.Range("I9").Value = 0.2967165721732315
$ws.Range("J9").Value = 0.2967165721732316
$ws.Range("O9").Value = 0.7176943460983047
$ws.Range("P9").Value = 0.7176943460983046
$ws.Range("S9").Value = 0.2129518062423978
$ws.Range("T9").Value = 0.2129518062423978

$ws.Range("I10").Value = 0.2967165721732315
$ws.Range("J10").Value = 0.2967165721732316
$ws.Range("O10").Value = 0.1936759918268401
$ws.Range("P10").Value = 0.1936759918268401
$ws.Range("S10").Value = 0.0574668764071108
$ws.Range("T10").Value = 0.05746687640711081
